$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$updates = @{
    "C2" = 39.19
    "F2" = 31.1
    "G2" = 0.04
    "H2" = 0.01
    "J2" = 0.07000000000000001
    "N2" = 0.03

    "C3" = 52.09
    "D3" = 34.2
    "I3" = 0.39
    "K3" = 0.13
    "M3" = 0.13
    "S3" = 0.13

    "C4" = 33.83
    "D4" = 15.18
    "G4" = 0.25
    "J4" = 0.16
    "K4" = 0.63
    "S4" = 0.03
    "T4" = 0.16

    "C5" = 28.8
    "H5" = 0.03
    "N5" = 0.03

    "C6" = 22.84
    "G6" = 0.03
    "J6" = 0.01
    "N6" = 0.07000000000000001
    "Q6" = 0.01
    "U6" = 0.03
    "Y6" = 0.02
    "AA6" = 0.05

    "C7" = 61.22
    "D7" = 48.83
    "F7" = 31.92
    "I7" = 0.15
    "J7" = 0.15
    "P7" = 0.15
    "Z7" = 0.29
    "AA7" = 0.29

    "C8" = 32.91
    "D8" = 40.98
    "I8" = 0.08
    "Z8" = 0.08

    "C9" = 28.98
    "G9" = 0.02
    "J9" = 0.04
    "S9" = 0
    "Y9" = 0.03

    "C10" = 55.62
    "D10" = 31.28
    "F10" = 7.58
    "G10" = 0.04
    "H10" = 0.03
    "I10" = 0.02
    "J10" = 0.01
    "S10" = 0.01
    "Y10" = 0.01

    "C11" = 34.77
    "G11" = 0.01
    "H11" = 0.01
    "J11" = 0.01
    "N11" = 0.02
    "Q11" = 0.03
    "U11" = 0.03
    "Y11" = 0.01

    "C12" = 23.45
    "D12" = 0.04
    "G12" = 0.06
    "J12" = 0.06
    "L12" = 0.1
    "S12" = 0.08
    "AA12" = 0.02

    "C13" = 31.32
    "F13" = 16.26
    "J13" = 0.03
    "S13" = 0.01
    "Y13" = 0.04

    "C14" = 31.12
    "D14" = 13.23
    "F14" = 37.57
    "G14" = 0.03
    "J14" = 0.01
    "Q14" = 0.01

    "C15" = 32.58
    "I15" = 0.03
    "J15" = 0.03
    "L15" = 0.01
    "N15" = 0.01
    "Y15" = 0.03

    "C16" = 38.58
    "D16" = 25.94
    "G16" = 0.04
    "L16" = 0.04
    "Q16" = 0.18
    "S16" = 0.04
    "Y16" = 0.04
    "AA16" = 0.11

    "C17" = 37.59
    "D17" = 39.49
    "H17" = 0.09
    "O17" = 0.63
    "T17" = 0.09
    "Z17" = 0.18
}

foreach ($addr in $updates.Keys) {
    $ws.Range($addr).Value = $updates[$addr]
}
